$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.046.96'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '1.854.81'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.93'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +6.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.329'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0694'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0990'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '2.123.44'
$ws.Range('E12').Value = '  +2.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.42'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.27%  '
$ws.Range('D14').Value = '1.856.88'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.677'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.68'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.92%  '
$ws.Range('D17').Value = '35.058.99'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.27'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').Value = '0.0₃0794'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.61'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.12'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.22%  '
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.06'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +27.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.93'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.63'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.124'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.40%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0556'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.99'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('E34').Value = '  +13.61%  '
$ws.Range('E35').Value = '  +23.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.28'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.777'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +12.90%  '
$ws.Range('E38').Value = '  +12.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '91.23'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0203'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.93%  '
$ws.Range('D41').Value = '1.350.40'
$ws.Range('E41').Value = '  +2.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.82'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.34'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.63'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +60.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.40'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0552'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +7.91%  '
$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.73'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.45'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.37%  '
$ws.Range('D49').Value = '2.037.60'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0679'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.23%  '
